$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.726.68"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.162.27"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.159.23"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "3.680.11"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "64.702.47"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "3.160.64"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  -6.38%  "
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.84%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").Value = "0.0₃0781"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "53.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").Value = "2.850.03"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("E45").Value = "  -3.77%  "
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.43%  "

Write-Output "Done"